# Generate Report for Handback
# Updates the handback-status workbook with freshly generated handoff/handback
# timestamps for the 24fb51b8-... file (zh-cn & de-de) and refreshes the
# "Latest HO Xliff Generate Date" for the ccd37122-... file on the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to ccd37122-6c86-4e63-b00d-1275dcf94100.md
$wsOverview.Range("G3").Value = "2016-08-22 22:47:27"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 2 corresponds to 24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md
$wsZhCn.Range("H2").Value = "2016-08-22 22:47:22"
$wsZhCn.Range("K2").Value = "2016-08-22 22:47:40"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 2 corresponds to 24fb51b8-bc83-4d13-ba2b-31a79a3a2db9.md
$wsDeDe.Range("H2").Value = "2016-08-22 22:47:27"
$wsDeDe.Range("K2").Value = "2016-08-22 22:47:47"
